# Add a new column C ('t+3') next to the existing A/B columns of predicted
# factors, mirroring the header formatting already used in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy B1's formatting (bold header / border / alignment style) onto C1,
# then overwrite the value with the new header number.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = 2

# Fill in the data rows for the new column.
$ws.Range("C2").Value = -4.693393667675708
$ws.Range("C3").Value = -1.171375753587524
$ws.Range("C4").Value = -0.06850756529651124
$ws.Range("C5").Value = -0.4260526063017495
$ws.Range("C6").Value = -0.1257255649142397
$ws.Range("C7").Value = 0.01878386401140213
